# Haver pull August 2023 — revise June/July 2023 rows and append new July 2023 row.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Single-cell value corrections to existing rows (revised Haver pull), plus the
# full rewrite of row 643 (2023-06-30) and a brand-new row 644 (2023-07-31).
$cellUpdates = @{
  "H132" = 935387
  "M132" = 14066934
  "M133" = 14469387
  "M134" = 14605022
  "M135" = 14690459
  "M136" = 14866918
  "M137" = 14861436
  "M138" = 14670828
  "M139" = 14467053
  "D140" = 823667
  "M140" = 14131833
  "N140" = 8935503
  "M141" = 13894129
  "N141" = 8782018
  "M142" = 13750036
  "N142" = 8776286
  "M143" = 13621905
  "N143" = 8836728
  "N144" = 9042536
  "N145" = 9395156
  "N146" = 9716020
  "N147" = 9916851
  "N148" = 10176484
  "N149" = 10402065
  "N150" = 10586465
  "N151" = 10834378
  "J156" = 37.94
  "D186" = 590412
  "G186" = 9817665
  "H186" = 1209170
  "I186" = 127.39
  "M186" = 13889959
  "N186" = 8211217
  "M187" = 13966192
  "N187" = 8231208
  "L188" = 125.89
  "M188" = 14177425
  "N188" = 8310994
  "M189" = 14260733
  "N189" = 8315354
  "M190" = 14423060
  "N190" = 8370935
  "D191" = 573557
  "G191" = 8670582
  "M191" = 14528237
  "N191" = 8368466
  "M192" = 14533012
  "N192" = 8290040
  "M193" = 14761535
  "N193" = 8365132
  "F194" = 14729052
  "M194" = 14844216
  "N194" = 8388019
  "M195" = 14889146
  "N195" = 8325491
  "F196" = 13097750
  "M196" = 14967828
  "N196" = 8319987
  "M197" = 15089972
  "N197" = 8341451
  "M198" = 15141528
  "N198" = 8330317
  "F199" = 10187285
  "M199" = 15309072
  "N199" = 8356565
  "E200" = 247048
  "M200" = 15440682
  "N200" = 8350584
  "O200" = 2599970
  "M201" = 15510456
  "N201" = 8368327
  "O201" = 2605545
  "F202" = 10071192
  "M202" = 15673766
  "N202" = 8381884
  "O202" = 2635460
  "F203" = 9851382
  "O203" = 2650960
  "O204" = 2658782
  "O205" = 2687331
  "O206" = 2682385
  "O207" = 2694736
  "O208" = 2719319
  "O209" = 2707761
  "O210" = 2680326
  "O211" = 2671972
  "B232" = 1361140
  "F232" = 11589035
  "R641" = 382073
  "P642" = 5216
  "Q642" = 14494
  "R642" = 386202
  "B643" = 1113048
  "C643" = 1654.7
  "D643" = 417231
  "E643" = 135748
  "F643" = 6939499
  "G643" = 5720156
  "H643" = 2439898
  "I643" = 434.29
  "J643" = 34.52
  "K643" = 13.95
  "L643" = 422.67
  "M643" = 27889394
  "N643" = 4820252
  "O643" = 1506166
  "P643" = 5246
  "Q643" = 14518
  "R643" = 387222

  "A644" = 45138
  "B644" = 1054535
  "D644" = 462890
  "E644" = 165809
  "F644" = 8551248
  "G644" = 6746674
  "H644" = 2803815
  "I644" = 423.2
  "J644" = 34.76
  "K644" = 14.1
  "L644" = 425.31
  "M644" = 28939334
  "N644" = 4915952
  "O644" = 1558501
  "P644" = 5235
  "Q644" = 14537

}

foreach ($ref in $cellUpdates.Keys) {
    $ws.Range($ref).Value = $cellUpdates[$ref]
}
